# Update covid_disparities_output workbook to the July 23, 2020 run results.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: Texas -- Bexar County ---
# Date advances a day; Total Cases / Total Deaths become unavailable (blank) in this run.
$ws.Range("B3").Value = 44035
$ws.Range("C3").ClearContents()
$ws.Range("D3").ClearContents()

# --- Row 10: Arkansas ---
$ws.Range("C10").Value = 36259
$ws.Range("D10").Value = 386
$ws.Range("E10").Value = 7710
$ws.Range("F10").Value = 102
$ws.Range("G10").Value = 24.54
$ws.Range("H10").Value = 26.36
$ws.Range("K10").Value = 31419
$ws.Range("L10").Value = 387

# --- Row 11: California - San Diego ---
$ws.Range("B11").Value = 44035
$ws.Range("C11").Value = 25608
$ws.Range("D11").Value = 512
$ws.Range("E11").Value = 942
$ws.Range("G11").Value = 4.7
$ws.Range("H11").Value = 3.81
$ws.Range("K11").Value = 20033
$ws.Range("L11").Value = 499

# --- Row 13: New Mexico ---
$ws.Range("B13").Value = 44035
$ws.Range("C13").Value = 18163
$ws.Range("D13").Value = 596
$ws.Range("E13").Value = 328

# --- Row 17: California - Los Angeles ---
$ws.Range("B17").Value = 44034
$ws.Range("C17").Value = 166848
$ws.Range("D17").Value = 4262
$ws.Range("E17").Value = 4451
$ws.Range("F17").Value = 424
$ws.Range("H17").Value = 10.66
$ws.Range("K17").Value = 95952
$ws.Range("L17").Value = 3976

# --- Row 27: Colorado ---
$ws.Range("B27").Value = 44035
$ws.Range("C27").Value = 42314
$ws.Range("D27").Value = 1786
$ws.Range("E27").Value = 2048
$ws.Range("F27").Value = 119
$ws.Range("G27").Value = 6.16
$ws.Range("H27").Value = 6.9
$ws.Range("K27").Value = 33225
$ws.Range("L27").Value = 1725

# --- Row 28: Nebraska ---
$ws.Range("B28").Value = 44035
$ws.Range("C28").Value = 23818
$ws.Range("D28").Value = 316
$ws.Range("E28").Value = 1414
$ws.Range("G28").Value = 7.69
$ws.Range("H28").Value = 7.72
$ws.Range("K28").Value = 18381
$ws.Range("L28").Value = 298

# --- Row 36: Washington ---
$ws.Range("B36").Value = 44035
$ws.Range("C36").Value = 50009
$ws.Range("D36").Value = 1482
$ws.Range("E36").Value = 1766
$ws.Range("G36").Value = 5.46
$ws.Range("H36").Value = 3.39
$ws.Range("K36").Value = 32323
$ws.Range("L36").Value = 1446

# --- Row 39: Delaware -- updated scrape error message ---
$ws.Range("O39").Value = "An error occurred. ... AttributeError(`"'numpy.float64' object has no attribute 'split'`")"
